$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 313 was missing its "symbol" (column B) value - fill it in to match
# the rest of the column.
$ws.Range("B313").Value = "ECONOMICS:CNCBBS"

# Append four more rows (314-317) that repeat the same datetime/OHLC/volume
# pattern as row 313, reusing row 313's date-cell formatting for column A.
# Row 317 mirrors the original row 313 shape (no symbol in column B).
$newRows = 314, 315, 316, 317
foreach ($r in $newRows) {
    $ws.Cells.Item(313, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = 45230

    if ($r -ne 317) {
        $ws.Cells.Item($r, 2).Value = "ECONOMICS:CNCBBS"
    }

    $ws.Cells.Item($r, 3).Value = 43325980000000
    $ws.Cells.Item($r, 4).Value = 43325980000000
    $ws.Cells.Item($r, 5).Value = 43325980000000
    $ws.Cells.Item($r, 6).Value = 43325980000000
    $ws.Cells.Item($r, 7).Value = 0
}

$excel.CutCopyMode = $false
